$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 19:22"

# --- Re-rank country name cells (country order shifted with the data refresh) ---
$ws.Range("A24").Value = "Irlanda"
$ws.Range("A25").Value = "India"
$ws.Range("A54").Value = "Sudafrica"
$ws.Range("A55").Value = "Argentina"
$ws.Range("A56").Value = "Grecia"
$ws.Range("A182").Value = "Belice"
$ws.Range("A183").Value = "Zimbabue"
$ws.Range("A184").Value = "Suazilandia"
$ws.Range("A185").Value = "Curazao"
$ws.Range("A186").Value = "Botsuana"
$ws.Range("A187").Value = "Malaui"
$ws.Range("A189").Value = "San Vicente y las Granadinas"
$ws.Range("A190").Value = "Nepal"

# --- Refresh numeric case-count data ---
$ws.Range("B4").Value = 545934
$ws.Range("C4").Value = 13055
$ws.Range("D4").Value = 31113
$ws.Range("E4").Value = 493347
$ws.Range("F4").Value = 11662
$ws.Range("G4").Value = 897
$ws.Range("H4").Value = 21474
$ws.Range("B7").Value = 132591
$ws.Range("C7").Value = 2937
$ws.Range("D7").Value = 27186
$ws.Range("E7").Value = 91012
$ws.Range("F7").Value = 6845
$ws.Range("G7").Value = 561
$ws.Range("H7").Value = 14393
$ws.Range("B8").Value = 126656
$ws.Range("C8").Value = 1204
$ws.Range("E8").Value = 66348
$ws.Range("G8").Value = 37
$ws.Range("H8").Value = 2908
$ws.Range("E15").Value = 12201
$ws.Range("G15").Value = 70
$ws.Range("H15").Value = 1106
$ws.Range("B16").Value = 23738
$ws.Range("C16").Value = 420
$ws.Range("E16").Value = 16403
$ws.Range("B17").Value = 21065
$ws.Range("C17").Value = 103
$ws.Range("E17").Value = 19748
$ws.Range("B24").Value = 9655
$ws.Range("C24").Value = 727
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 9296
$ws.Range("F24").Value = 194
$ws.Range("G24").Value = 14
$ws.Range("H24").Value = 334
$ws.Range("B25").Value = 9166
$ws.Range("C25").Value = 720
$ws.Range("D25").Value = 1061
$ws.Range("E25").Value = 7780
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 325
$ws.Range("E33").Value = 5132
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = 316
$ws.Range("B54").Value = 2173
$ws.Range("C54").Value = 145
$ws.Range("D54").Value = 410
$ws.Range("E54").Value = 1738
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 25
$ws.Range("B55").Value = 2142
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 468
$ws.Range("E55").Value = 1584
$ws.Range("F55").Value = 83
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 90
$ws.Range("B56").Value = 2114
$ws.Range("C56").Value = 33
$ws.Range("D56").Value = 269
$ws.Range("E56").Value = 1747
$ws.Range("F56").Value = 76
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 98
$ws.Range("B117").Value = 210
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 56
$ws.Range("E117").Value = 147
$ws.Range("C182").Value = 1
$ws.Range("E182").Value = 12
$ws.Range("F182").Value = 1
$ws.Range("H182").Value = 2
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 11
$ws.Range("H183").Value = 3
$ws.Range("C184").Value = 2
$ws.Range("E184").Value = 7
$ws.Range("H184").Value = 0
$ws.Range("B185").Value = 14
$ws.Range("D185").Value = 7
$ws.Range("E185").Value = 6
$ws.Range("C186").Value = 0
$ws.Range("E186").Value = 12
$ws.Range("F186").Value = 0
$ws.Range("H186").Value = 1
$ws.Range("C187").Value = 1
$ws.Range("C189").Value = 0
$ws.Range("C190").Value = 3
